$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"0.0292345"
$ws.Range("H2").Value = [double]"0.058469"
$ws.Range("I2").Value = [double]"0.4428765120700495"
$ws.Range("J2").Value = [double]"0.346386487911515"
$ws.Range("M2").Value = [double]"0.2272265"
$ws.Range("N2").Value = [double]"0.454453"
$ws.Range("O2").Value = [double]"0.08704083604617911"
$ws.Range("P2").Value = [double]"0.08229687998280369"
$ws.Range("Q2").Value = [double]"0.00664285311425"
$ws.Range("R2").Value = [double]"0.026571412457"
$ws.Range("S2").Value = [double]"0.03854834187579284"
$ws.Range("T2").Value = [double]"0.02850652722331883"
$ws.Range("G3").Value = [double]"0.0292345"
$ws.Range("H3").Value = [double]"0.058469"
$ws.Range("I3").Value = [double]"0.4428765120700495"
$ws.Range("J3").Value = [double]"0.346386487911515"
$ws.Range("O3").Value = [double]"0.1052353694185077"
$ws.Range("P3").Value = [double]"0.149249644656207"
$ws.Range("Q3").Value = [double]"0.008031438267666667"
$ws.Range("R3").Value = [double]"0.048188629606"
$ws.Range("S3").Value = [double]"0.04660627335447184"
$ws.Range("T3").Value = [double]"0.05169806023450516"
$ws.Range("G4").Value = [double]"0.0292345"
$ws.Range("H4").Value = [double]"0.058469"
$ws.Range("I4").Value = [double]"0.4428765120700495"
$ws.Range("J4").Value = [double]"0.346386487911515"
$ws.Range("M4").Value = [double]"2.082377"
$ws.Range("N4").Value = [double]"4.164754"
$ws.Range("O4").Value = [double]"0.7976703203338269"
$ws.Range("P4").Value = [double]"0.7541951755096822"
$ws.Range("Q4").Value = [double]"0.0608772504065"
$ws.Range("R4").Value = [double]"0.243509001626"
$ws.Range("S4").Value = [double]"0.3532694492512443"
$ws.Range("T4").Value = [double]"0.2612430180446075"
$ws.Range("G5").Value = [double]"0.0292345"
$ws.Range("H5").Value = [double]"0.058469"
$ws.Range("I5").Value = [double]"0.4428765120700495"
$ws.Range("J5").Value = [double]"0.346386487911515"
$ws.Range("M5").Value = [double]"0.02610733333333333"
$ws.Range("N5").Value = [double]"0.078322"
$ws.Range("O5").Value = [double]"0.0100006122537187"
$ws.Range("P5").Value = [double]"0.01418332860386696"
$ws.Range("Q5").Value = [double]"0.0007632348363333333"
$ws.Range("R5").Value = [double]"0.004579409018"
$ws.Range("S5").Value = [double]"0.004429036273491936"
$ws.Range("T5").Value = [double]"0.004912913381988407"
$ws.Range("G6").Value = [double]"0.0292345"
$ws.Range("H6").Value = [double]"0.058469"
$ws.Range("I6").Value = [double]"0.4428765120700495"
$ws.Range("J6").Value = [double]"0.346386487911515"
$ws.Range("M6").Value = [double]"0.000138"
$ws.Range("N6").Value = [double]"0.000414"
$ws.Range("O6").Value = [double]"5.286194776741585E-05"
$ws.Range("P6").Value = [double]"7.49712474400669E-05"
$ws.Range("Q6").Value = [double]"4.034361E-06"
$ws.Range("R6").Value = [double]"2.4206166E-05"
$ws.Range("S6").Value = [double]"2.341131504846227E-05"
$ws.Range("T6").Value = [double]"2.596902709510993E-05"
$ws.Range("I7").Value = [double]"0.5571234879299505"
$ws.Range("J7").Value = [double]"0.6536135120884849"
$ws.Range("M7").Value = [double]"0.2272265"
$ws.Range("N7").Value = [double]"0.454453"
$ws.Range("O7").Value = [double]"0.08704083604617911"
$ws.Range("P7").Value = [double]"0.08229687998280369"
$ws.Range("Q7").Value = [double]"0.008356481763999999"
$ws.Range("R7").Value = [double]"0.050138890584"
$ws.Range("S7").Value = [double]"0.04849249417038626"
$ws.Range("T7").Value = [double]"0.05379035275948486"
$ws.Range("I8").Value = [double]"0.5571234879299505"
$ws.Range("J8").Value = [double]"0.6536135120884849"
$ws.Range("O8").Value = [double]"0.1052353694185077"
$ws.Range("P8").Value = [double]"0.149249644656207"
$ws.Range("S8").Value = [double]"0.05862909606403585"
$ws.Range("T8").Value = [double]"0.09755158442170185"
$ws.Range("I9").Value = [double]"0.5571234879299505"
$ws.Range("J9").Value = [double]"0.6536135120884849"
$ws.Range("M9").Value = [double]"2.082377"
$ws.Range("N9").Value = [double]"4.164754"
$ws.Range("O9").Value = [double]"0.7976703203338269"
$ws.Range("P9").Value = [double]"0.7541951755096822"
$ws.Range("Q9").Value = [double]"0.076581496552"
$ws.Range("R9").Value = [double]"0.459488979312"
$ws.Range("S9").Value = [double]"0.4444008710825825"
$ws.Range("T9").Value = [double]"0.4929521574650747"
$ws.Range("I10").Value = [double]"0.5571234879299505"
$ws.Range("J10").Value = [double]"0.6536135120884849"
$ws.Range("M10").Value = [double]"0.02610733333333333"
$ws.Range("N10").Value = [double]"0.078322"
$ws.Range("O10").Value = [double]"0.0100006122537187"
$ws.Range("P10").Value = [double]"0.01418332860386696"
$ws.Range("Q10").Value = [double]"0.0009601232906666665"
$ws.Range("R10").Value = [double]"0.008641109616"
$ws.Range("S10").Value = [double]"0.005571575980226767"
$ws.Range("T10").Value = [double]"0.009270415221878551"
$ws.Range("I11").Value = [double]"0.5571234879299505"
$ws.Range("J11").Value = [double]"0.6536135120884849"
$ws.Range("M11").Value = [double]"0.000138"
$ws.Range("N11").Value = [double]"0.000414"
$ws.Range("O11").Value = [double]"5.286194776741585E-05"
$ws.Range("P11").Value = [double]"7.49712474400669E-05"
$ws.Range("Q11").Value = [double]"5.075088E-06"
$ws.Range("R11").Value = [double]"4.5675792E-05"
$ws.Range("S11").Value = [double]"2.945063271895358E-05"
$ws.Range("T11").Value = [double]"4.900222034495696E-05"
